# Adapt column header formatting to respective input file names (#7)
#
# - Rename the "_old"/"_new" header-name suffixes to "_FV2404"/"_FV2410"
# - Turn the data range A1:U84 into an Excel Table ("Table1") with an AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells in row 1 (A1:U1) ------------------------------
# "<name>_old" -> "<name>_FV2404", "<name>_new" -> "<name>_FV2410"
$lastCol = 21
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2
    if ($current -ne $null) {
        $updated = $current -replace "_old$", "_FV2404"
        $updated = $updated -replace "_new$", "_FV2410"
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}

# --- 2. Turn A1:U84 into a real Table (adds AutoFilter + table part) ------
$dataRange = $ws.Range("A1:U84")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
